$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Consolidate the title's "A slide" text runs into a single run.
# (Setting through a temporary distinct value forces the run-merge logic
# to actually rewrite the paragraph instead of treating it as a no-op
# when the concatenated text already reads the same.)
$titleShape = $s.Shapes.Item("Title 1")
$titleShape.TextFrame.TextRange.Text = "__tmp__"
$titleShape.TextFrame.TextRange.Text = "A slide"

# Consolidate the "Just an image on this side" textbox runs into a single run.
$captionShape = $s.Shapes.Item("TextBox 3")
$captionShape.TextFrame.TextRange.Text = "__tmp__"
$captionShape.TextFrame.TextRange.Text = "Just an image on this side"
